$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.874.77"
$ws.Range("E2").Value = "  -1.00%  "

# Row 3
$ws.Range("D3").Value = "2.356.82"
$ws.Range("E3").Value = "  -4.43%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "542.35"
$ws.Range("E5").Value = "  -1.63%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.59"
$ws.Range("E6").Value = "  -6.41%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.525"
$ws.Range("E8").Value = "  -11.56%  "

# Row 9
$ws.Range("D9").Value = "2.355.74"
$ws.Range("E9").Value = "  -4.54%  "

# Row 10
$ws.Range("E10").Value = "  -2.79%  "

# Row 11
$ws.Range("E11").Value = "  +0.11%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.25"
$ws.Range("E12").Value = "  -3.80%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.341"
$ws.Range("E13").Value = "  -3.70%  "

# Row 14
$ws.Range("E14").Value = "  -5.12%  "

# Row 15
$ws.Range("D15").Value = "2.776.06"
$ws.Range("E15").Value = "  -4.47%  "

# Row 16
$ws.Range("D16").Value = "60.643.51"
$ws.Range("E16").Value = "  -1.22%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000161"
$ws.Range("E17").Value = "  -3.77%  "

# Row 18
$ws.Range("D18").Value = "2.356.59"
$ws.Range("E18").Value = "  -4.27%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.61"
$ws.Range("E19").Value = "  -4.79%  "

# Row 20
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.06"
$ws.Range("E20").Value = "  -3.22%  "

# Row 21
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "314.32"
$ws.Range("E21").Value = "  -1.41%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.58"
$ws.Range("E22").Value = "  -7.89%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.21%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.87"
$ws.Range("E24").Value = "  -1.22%  "

# Row 25
$ws.Range("E25").Value = "  -1.22%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.31"
$ws.Range("E26").Value = "  +7.83%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.09%  "

# Row 28
$ws.Range("D28").Value = "2.468.73"
$ws.Range("E28").Value = "  -4.34%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0894"
$ws.Range("E29").Value = "  -9.29%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.91"
$ws.Range("E30").Value = "  -4.94%  "

# Row 31
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "502.18"
$ws.Range("E31").Value = "  -9.30%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.39"
$ws.Range("E32").Value = "  -6.59%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.145"
$ws.Range("E33").Value = "  -1.63%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.78"
$ws.Range("E34").Value = "  -6.51%  "

# Row 35
$ws.Range("E35").Value = "  -3.57%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  +0.02%  "

# Row 37
$ws.Range("E37").Value = "  -6.44%  "

# Row 38
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.50"
$ws.Range("E38").Value = "  +0.29%  "

# Row 39
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.372"
$ws.Range("E39").Value = "  -2.12%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.21"
$ws.Range("E40").Value = "  -11.91%  "

# Row 41
$ws.Range("E41").Value = "  +1.12%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "138.91"
$ws.Range("E42").Value = "  -2.48%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.07"
$ws.Range("E44").Value = "  -1.01%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.15"
$ws.Range("E45").Value = "  -9.80%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "138.95"
$ws.Range("E46").Value = "  -5.24%  "

# Row 47
$ws.Range("E47").Value = "  -2.61%  "

# Row 48
$ws.Range("E48").Value = "  -5.22%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.50"
$ws.Range("E49").Value = "  -9.85%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.570"
$ws.Range("E50").Value = "  -3.72%  "

# Row 51
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0223"
$ws.Range("E51").Value = "  -2.84%  "
